# Update classification-report metrics (precision/recall/f1-score columns)
# for rows 2-26 with newly computed values from the latest model run
# (TPE search over decision tree / LightGBM classifiers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5952380952380952
$ws.Range("C2").Value = 0.6578947368421053
$ws.Range("D2").Value = 0.625
$ws.Range("B3").Value = 0.7450980392156863
$ws.Range("C3").Value = 0.6909090909090909
$ws.Range("D3").Value = 0.7169811320754716
$ws.Range("B4").Value = 0.6774193548387096
$ws.Range("C4").Value = 0.6774193548387096
$ws.Range("D4").Value = 0.6774193548387096
$ws.Range("E4").Value = 0.6774193548387096
$ws.Range("B5").Value = 0.6701680672268908
$ws.Range("C5").Value = 0.6744019138755981
$ws.Range("D5").Value = 0.6709905660377358
$ws.Range("B6").Value = 0.6838649438269933
$ws.Range("C6").Value = 0.6774193548387096
$ws.Range("D6").Value = 0.6793974437005478
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 0.1052631578947368
$ws.Range("D7").Value = 0.1904761904761905
$ws.Range("B8").Value = 0.6179775280898876
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0.7638888888888888
$ws.Range("B9").Value = 0.6344086021505376
$ws.Range("C9").Value = 0.6344086021505376
$ws.Range("D9").Value = 0.6344086021505376
$ws.Range("E9").Value = 0.6344086021505376
$ws.Range("B10").Value = 0.8089887640449438
$ws.Range("C10").Value = 0.5526315789473684
$ws.Range("D10").Value = 0.4771825396825397
$ws.Range("B11").Value = 0.7740727316660626
$ws.Range("C11").Value = 0.6344086021505376
$ws.Range("D11").Value = 0.5295912271718723
$ws.Range("B12").Value = 0.5714285714285714
$ws.Range("C12").Value = 0.7368421052631579
$ws.Range("D12").Value = 0.6436781609195403
$ws.Range("B13").Value = 0.7727272727272727
$ws.Range("C13").Value = 0.6181818181818182
$ws.Range("D13").Value = 0.6868686868686869
$ws.Range("B14").Value = 0.6666666666666666
$ws.Range("C14").Value = 0.6666666666666666
$ws.Range("D14").Value = 0.6666666666666666
$ws.Range("E14").Value = 0.6666666666666666
$ws.Range("B15").Value = 0.6720779220779221
$ws.Range("C15").Value = 0.6775119617224881
$ws.Range("D15").Value = 0.6652734238941136
$ws.Range("B16").Value = 0.6904761904761905
$ws.Range("C16").Value = 0.6666666666666666
$ws.Range("D16").Value = 0.6692209450830141
$ws.Range("B17").Value = 0.576271186440678
$ws.Range("C17").Value = 0.8947368421052632
$ws.Range("D17").Value = 0.7010309278350517
$ws.Range("B18").Value = 0.8823529411764706
$ws.Range("C18").Value = 0.5454545454545454
$ws.Range("D18").Value = 0.6741573033707865
$ws.Range("B19").Value = 0.6881720430107527
$ws.Range("C19").Value = 0.6881720430107527
$ws.Range("D19").Value = 0.6881720430107527
$ws.Range("E19").Value = 0.6881720430107527
$ws.Range("B20").Value = 0.7293120638085743
$ws.Range("C20").Value = 0.7200956937799043
$ws.Range("D20").Value = 0.6875941156029191
$ws.Range("B21").Value = 0.7572872779510931
$ws.Range("C21").Value = 0.6881720430107527
$ws.Range("D21").Value = 0.685137924119626
$ws.Range("B22").Value = 0.5918367346938775
$ws.Range("C22").Value = 0.7631578947368421
$ws.Range("D22").Value = 0.6666666666666667
$ws.Range("B23").Value = 0.7954545454545454
$ws.Range("C23").Value = 0.6363636363636364
$ws.Range("D23").Value = 0.7070707070707071
$ws.Range("B24").Value = 0.6881720430107527
$ws.Range("C24").Value = 0.6881720430107527
$ws.Range("D24").Value = 0.6881720430107527
$ws.Range("E24").Value = 0.6881720430107527
$ws.Range("B25").Value = 0.6936456400742115
$ws.Range("C25").Value = 0.6997607655502392
$ws.Range("D25").Value = 0.6868686868686869
$ws.Range("B26").Value = 0.7122558700899715
$ws.Range("C26").Value = 0.6881720430107527
$ws.Range("D26").Value = 0.6905615292712067
